$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- First-page header: BTec logo "image1.jpg" -> "image2.jpg" ---
$hdrFirst = $sec.Headers.Item(2)
$hdrShape = $hdrFirst.Range.InlineShapes.Item(1)
$hdrShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

# --- Default (primary) footer: Pearson logo "image2.png" -> "image1.png" ---
$ftrDefault = $sec.Footers.Item(1)
$ftrDefaultShape = $ftrDefault.Range.InlineShapes.Item(1)
$ftrDefaultShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# --- First-page footer: Pearson logo "image2.png" -> "image1.png" ---
$ftrFirst = $sec.Footers.Item(2)
$ftrFirstShape = $ftrFirst.Range.InlineShapes.Item(1)
$ftrFirstShape.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

Write-Output "Renamed header/footer picture names."
